$wb = $excel.ActiveWorkbook

$alpha = $wb.Worksheets.Item("alphaSheet")
$beta  = $wb.Worksheets.Item("betaSheet")

# --- Add the new "realTime" worksheet as the last tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$rt = $wb.Worksheets.Add($null, $lastSheet)
$rt.Name = "realTime"

# Reuse the existing (already-present) cell formats from betaSheet so no
# new style entries are introduced - betaSheet already has the highlighted
# header style (A3-like) and the plain bordered styles (B3-like, A2-like).
$beta.Range("A3").Copy()
$rt.Range("A1").PasteSpecial(-4122)
$rt.Range("D1").PasteSpecial(-4122)

$beta.Range("B3").Copy()
$rt.Range("B1").PasteSpecial(-4122)
$rt.Range("C1").PasteSpecial(-4122)
$rt.Range("B2").PasteSpecial(-4122)

$beta.Range("A2").Copy()
$rt.Range("A2").PasteSpecial(-4122)

$beta.Range("B2").Copy()
$rt.Range("C2").PasteSpecial(-4122)

# Populate the new sheet (order matters for shared-string table ordering)
$rt.Range("A1").Value = "HighLightElement"
$rt.Range("B1").Value = "userName"
$rt.Range("C1").Value = "passWord"
$rt.Range("D1").Value = "HighLightElement"
$rt.Range("B2").Value = "Admin"
$rt.Range("C2").Value = "admin123"

# --- Update existing data on alphaSheet ---
$alpha.Range("D10").Value = "12456"
$alpha.Range("D6").Value = "786"

# --- Selections (restore per-sheet selection state) ---
$alpha.Range("D6").Select()
$beta.Range("A3:F4").Select()
$rt.Range("C2").Select()
